# Correction de nom de variable : "T_compo_carte" -> "T_compo_paquet"
#
# The bold run that used to read "_carte" must become two runs reading
# "_" and "paquet" (both bold), with Word's automatic "_GoBack" bookmark
# ending up right between them (because that's where the user's last edit
# landed). The bookmark therefore has to disappear from its old location
# further down in the document (between "cartes" and "[compteur].").

$d = $word.ActiveDocument

# Locate the run that currently holds "_carte".
$target = $d.Content
$found = $target.Find.Execute("_carte", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find '_carte' in the document"
}
$fullStart = $target.Start
$fullEnd = $target.End

# Character boundaries:
#   $boundary1 -> between the previous run ("_compo") and "_carte"
#   $boundary2 -> between "_" and "carte" inside the "_carte" run
$boundary1 = $fullStart
$boundary2 = $fullStart + 1

# Drop a throw-away bookmark on $boundary1 and remove it immediately: Word
# always splits the run at the insertion point, and - unlike a plain text
# edit - that split survives even after the bookmark itself is deleted.
# This stops the upcoming text edit from merging back into the identically
# formatted "_compo" run next door.
$barrier = $d.Range($boundary1, $boundary1)
$d.Bookmarks.Add("ZZ_tmp_barrier", $barrier) | Out-Null
$d.Bookmarks.Item("ZZ_tmp_barrier").Delete()

# Plant the real "_GoBack" bookmark at $boundary2 - this both marks the
# final cursor position (moving it off its old spot automatically, since a
# document can only have one "_GoBack") and acts as a second barrier so the
# "carte" -> "paquet" replacement below cannot merge "_" and "paquet" back
# into a single run.
$goBack = $d.Range($boundary2, $boundary2)
$d.Bookmarks.Add("_GoBack", $goBack) | Out-Null

# Replace "carte" with "paquet", keeping the bold formatting of the run.
$word.ActiveDocument.Range($boundary2, $fullEnd).Text = "paquet"
